$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new values on existing row 2 (PriceChange / UpDown columns)
$ws.Range("X2").Value = -0.059997999999993112
$ws.Range("Y2").Value = "Down"

# Add new row 3 with a fresh scan result
$ws.Range("A3").Value = 42648.663935185185
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = "Buy"
$ws.Range("D3").Value = 38
$ws.Range("E3").Value = 8602
$ws.Range("F3").Value = 1007
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 38
$ws.Range("I3").Value = 84
$ws.Range("J3").Value = 15
$ws.Range("K3").Value = 6117
$ws.Range("L3").Value = 137
$ws.Range("M3").Value = 88
$ws.Range("N3").Value = 50
$ws.Range("O3").Value = 9
$ws.Range("P3").Value = "Bag"
$ws.Range("Q3").Value = 48.098617091043238
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.096699999999999994
$ws.Range("T3").Value = 0.027400000000000001
$ws.Range("U3").Value = 4.8
$ws.Range("V3").Value = 2.2799999999999998
$ws.Range("W3").Value = 0

# Copy number formats (date / percentage) from row 2 without introducing
# new style entries, matching the original workbook's style table.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)

$ws.Range("S2").Copy()
$ws.Range("S3").PasteSpecial(-4122)

$ws.Range("T2").Copy()
$ws.Range("T3").PasteSpecial(-4122)
